# The workbook has two sheets: "Mondai" (the blank question template the
# user fills in) and "作成例" (a worked example kept for reference).
#
# This edit:
#   1. Clears out the two sample rows of data (row 2 and row 3) on the
#      "Mondai" sheet, leaving the header row and cell formatting intact,
#      so the template starts empty again.
#   2. Tweaks the guidance text in the "作成例" sheet's notes section
#      (cell B14) to clarify that the star/plus value "1" must be typed
#      using half-width (ASCII) digits.

$wb = $excel.ActiveWorkbook

$mondai = $wb.Worksheets.Item("Mondai")
$mondai.Range("A2:L3").ClearContents()

$sakuseirei = $wb.Worksheets.Item("作成例")
$sakuseirei.Range("B14").Value = "スター＋（必須）：「1」を半角数字で入力"
